# Add two new slides (Title and Content layout) describing the
# definition and classification of software requirements.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 4 : "Định nghĩa" (Definition)
# ---------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)

$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Định nghĩa"
$s4.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-US"

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "Yêu cầu cho 1 phần mềm cụ thể là tổng hợp những yêu cầu từ nhiều người khác nhau về tổ chức, mức độ chuyên môn và mức độ tham gia, tương tác với phần mềm trong môi trường hoạt động của nó.`r`rCó thể kiểm chứng một cách riêng rẽ ở mức chức năng(yêu cầu chức năng) hoặc mức hệ thống(yêu cầu phi chức năng)`r`rCung cấp các chỉ số đánh giá độ ưu tiên về các mặt khi cân nhắc về nguồn tài nguyên.`r`rCung cấp các giá trị trạng thái để theo dõi tiến độ của dự án."
$body4.Font.Size = 18
$body4.LanguageID = "vi-VN"

# ---------------------------------------------------------------
# Slide 5 : "Phân loại" (Classification)
# ---------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)

$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Phân loại"
$s5.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-US"

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Theo sản phẩm và tiến trình`rYêu cầu sản phẩm: là những đòi hỏi hay ràng buộc mà phần mềm phải thực hiện.`r`rYêu cầu tiến trình: là những ràng buộc liên quan đến việc phát triển phần mềm đó(quy trình, đối tác kiểm thử, phân tích, kĩ thuật sử dụng,...).`rTheo chức năng`rYêu cầu chức năng: đặc tả các chức năng mà phần mềm phải thực hiện.`r`rYêu cầu phi chức năng: là các ràng buộc về giải pháp và chất lượng(hiệu năng, việc bảo trì, độ an toàn, bảo mật,...)."
$body5.LanguageID = "vi-VN"

# Paragraph 1: "Theo sản phẩm và tiến trình" - level 0, size 20, bold
$para = $body5.Paragraphs(1)
$para.Font.Size = 20
$para.Font.Bold = $true

# Paragraphs 2-4: product/process requirement bullets - level 1, size 16
for ($i = 2; $i -le 4; $i++) {
    $para = $body5.Paragraphs($i)
    $para.IndentLevel = 2
    $para.Font.Size = 16
    $para.ParagraphFormat.Bullet.Font.Name = "Wingdings"
    $para.ParagraphFormat.Bullet.Character = 216
}

# Paragraph 5: "Theo chức năng" - level 0, size 20, bold
$para = $body5.Paragraphs(5)
$para.Font.Size = 20
$para.Font.Bold = $true

# Paragraphs 6-8: functional/non-functional requirement bullets - level 1, size 16
for ($i = 6; $i -le 8; $i++) {
    $para = $body5.Paragraphs($i)
    $para.IndentLevel = 2
    $para.Font.Size = 16
    $para.ParagraphFormat.Bullet.Font.Name = "Wingdings"
    $para.ParagraphFormat.Bullet.Character = 216
}

Write-Output "Slides added. Total slide count: $($p.Slides.Count)"
